# Adds a new "2020" column (L) to the 16.7.1.1 table, mirroring the
# formatting of the existing "2019" column (K), then updates the sheet's
# view/selection to the newly-added range (L4:L13), scrolled so column C
# is left-most.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header + data values for the new column L (rows 4-13).
$values = @{
    4  = 2020
    5  = 1.2
    6  = 1.7
    7  = 0.4
    8  = 3.3
    9  = 3.9
    10 = 2.4
    11 = 95.5
    12 = 94.4
    13 = 97.2
}

foreach ($row in $values.Keys | Sort-Object) {
    $ws.Cells.Item($row, 12).Value = $values[$row]   # column 12 = L
}

# Match column K's cell formatting (number format/font/border/etc.) for
# each new column L cell so the look is identical to the rest of the row.
foreach ($row in $values.Keys | Sort-Object) {
    [void]$ws.Cells.Item($row, 11).Copy()            # column 11 = K
    [void]$ws.Cells.Item($row, 12).PasteSpecial(-4122)   # xlPasteFormats
}
$excel.CutCopyMode = $false

# Scroll the view so column C is the left-most visible column, and select
# the newly added L4:L13 range (active cell L4).
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
[void]$ws.Range("L4:L13").Select()
